# BomPlan.xlsx update
# 1. Project No and Name complete
# 2. Description for SA complete
# 3. Color of SA in BOM Plan and of Detail Design in Fixed weight complete

$wb = $excel.ActiveWorkbook

# --- Sheet "BOMPlanTest": Detail Design fixed-weight color codes (col I) ---
$ws1 = $wb.Worksheets.Item("BOMPlanTest")
$ws1.Range("I4").Value = "RAL5010"
$ws1.Range("I5").Value = "RAL5011"
$ws1.Range("I6").Value = "RAL5012"
$ws1.Range("I7").Value = "RAL5012"
$ws1.Range("I8").Value = "RAL5013"
$ws1.Range("I9").Value = "RAL5013"

# --- Sheet "FAB01-B3-02": BOM Plan for SA ---
$ws2 = $wb.Worksheets.Item("FAB01-B3-02")

# Project No and Name
$ws2.Range("B1").Value = "S32A1305700"

# Column header
$ws2.Range("I3").Value = "Color Code"

# Description for SA
$ws2.Range("B4").Value = "SA-Test1"
$ws2.Range("B5").Value = "SA-Test1"
$ws2.Range("B6").Value = "SA-Test2"
$ws2.Range("B8").Value = "SA-Test3"

# Type code
$ws2.Range("G8").Value = "TS"

# Color of SA in BOM Plan (col I)
$ws2.Range("I4").Value = "RBL5010"
$ws2.Range("I5").Value = "RBL5011"
$ws2.Range("I6").Value = "RBL5012"
$ws2.Range("I7").Value = "RBL5012"
$ws2.Range("I8").Value = "RBL5013"
$ws2.Range("I9").Value = "RBL5013"
